# Auto-update draw results: append the 2025-10-18 "Pick 4" draw as a new
# row at the bottom of the "Results" sheet, mirroring the existing rows
# (Date, Game, Phase, Result, InsertedAt), all stored as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 32

$date       = "2025-10-18"
$game       = "Pick 4"
$phase      = "251018"
$result     = "2-8-1-4"
$insertedAt = "2025-10-18T21:35:23.528+04:00"

# Columns A (date-shaped) and C (digits-only) would otherwise be
# auto-converted to a date serial / a number by the smart Value parser,
# like every other row in this sheet they must stay literal text.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("C$newRow").NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $date
$ws.Cells.Item($newRow, 2).Value = $game
$ws.Cells.Item($newRow, 3).Value = $phase
$ws.Cells.Item($newRow, 4).Value = $result
$ws.Cells.Item($newRow, 5).Value = $insertedAt

# Drop the explicit "text" number format again now that the values are
# committed as text, so the new cells don't end up with a style index
# that the rest of the (unstyled) sheet doesn't carry.
$ws.Range("A$newRow").ClearFormats()
$ws.Range("C$newRow").ClearFormats()
